$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D/E columns hold text-formatted price/volume strings (e.g. "26.175.54", "  -4.52%  ").
# Force Text number format first so Excel does not auto-coerce these into numeric values.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.175.54'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -4.52%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.658.27'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -3.11%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.32%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.03'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -2.85%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5153'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -3.47%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.007'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.19%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2566'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -3.90%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06395'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -3.28%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.82'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -5.25%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07786'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.95%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.672.81'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.94%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.303'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -5.61%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.884.52'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -3.20%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -4.09%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0₅8029'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.62%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.19'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -5.46%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '26.201.86'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -4.37%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.32%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '210.13'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -2.80%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.393'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -5.83%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -3.81%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.872'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -1.77%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.008'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.22%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '143.72'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.81%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.762'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +1.85%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1161'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -4.58%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.960'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -4.44%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -3.09%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05247'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -2.80%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -2.67%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.372'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -3.38%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.207'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -6.40%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -4.82%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.750'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -4.57%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.364'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -2.23%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9238'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5716'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -1.99%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.155.32'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +10.63%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01589'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.69%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.008'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.22%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8379'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.39%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.662'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -3.42%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '99.81'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.05%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.794.87'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -3.24%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0₈110'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -7.40%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4505'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.32%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '56.02'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -3.42%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.53%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.903'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -2.24%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05091'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -2.89%  '
